$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 4 (Cruz / Female / 29125) -------------------------------------
# Build each cell's formatting by copying the matching cell directly above it
# (same column, row 3) so the new row reuses the existing style indices
# instead of minting new ones.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").ClearContents()

$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 29125

$ws.Range("D3").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# --- New column E (Voting behavior / Last date voted / 11/10/2012) ---------
# E1 starts out as a copy of the "Demographics" banner cell (fill + border +
# centered) then drops the centering, matching the new un-merged header cell.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").HorizontalAlignment = 1

$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# --- Cell text/values (order controls shared-string table order) -----------
$ws.Range("E2").Value = "Last date voted"
$ws.Range("E1").Value = "Voting behavior"
$ws.Range("B4").Value = "Cruz"
$ws.Range("C4").Value = "Female"
$ws.Range("E3").Value = 41223

# --- Header banner fill: theme Accent5 -> Accent1 (same tint) --------------
# (MsoThemeColorIndex 5 = Accent1, matching OOXML fgColor theme="4")
$ws.Range("E1").Interior.ThemeColor = 5
$ws.Range("A1:D1").Interior.ThemeColor = 5

# --- Column widths -----------------------------------------------------
$ws.Columns("C").ColumnWidth = 6.6
$ws.Columns("E").ColumnWidth = 18.6

# --- Selection -----------------------------------------------------------
$ws.Range("H6").Select()
